# The workbook's only sheet ("Sheet1") carries protection on its data cells,
# so unprotect first, make the data edits, then restore protection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure text (A11)
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-08 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-8
$ws.Range("D2").Value = 0.5300987286341232
$ws.Range("E2").Value = -0.009481240688067105

$ws.Range("D3").Value = 0.2694323540539976
$ws.Range("E3").Value = -0.007098848111438572

$ws.Range("D4").Value = 0.05017291184907132
$ws.Range("E4").Value = -0.01174434417109649

$ws.Range("D5").Value = 0.09406314369699235
$ws.Range("E5").Value = -0.01218108574717069

$ws.Range("D6").Value = 0.02684827908587062
$ws.Range("E6").Value = -0.01244731941585808

$ws.Range("D7").Value = 0.02938458267994476
$ws.Range("E7").Value = -0.01174393613754998

$ws.Range("E8").Value = -0.009352971924168285

# Restore sheet protection (present in the original workbook)
$ws.Protect()
